$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows for Mitutoyo scales re-measured on the Canon EOS R6 Mark II
$ws.Range("A93").Value = "Canon EOS R6 Mark II"
$ws.Range("B93").Value = "Mitutoyo"
$ws.Range("C93").Value = "Mitutoyo HR 5.0x"
$ws.Range("D93").Value = "5.0x"
$ws.Range("E93").Value = "6000"
$ws.Range("F93").Value = "4000"
$ws.Range("G93").Value = "892"
$ws.Range("H93").Value = "1000"
$ws.Range("I93").Value = "µm"

$ws.Range("A94").Value = "Canon EOS R6 Mark II"
$ws.Range("B94").Value = "Mitutoyo"
$ws.Range("C94").Value = "Mitutoyo 7.5x"
$ws.Range("D94").Value = "7.5x"
$ws.Range("E94").Value = "6000"
$ws.Range("F94").Value = "4000"
$ws.Range("G94").Value = "673"
$ws.Range("H94").Value = "500"
$ws.Range("I94").Value = "µm"

# Copy formatting from the preceding analogous block (rows 91-92, the other
# Mitutoyo / R5 entries) so the new rows look consistent with the rest of
# the table.
$ws.Range("A91:I92").Copy()
$ws.Range("A93:I94").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update view to match the saved state after the edit
$ws.Range("G95").Select()
$excel.ActiveWindow.ScrollRow = 68
$excel.ActiveWindow.ScrollColumn = 1

$wb.Save()
